# Updated cryptos list on Fri May 17 17:27:12 UTC 2024 with GitHub Actions
#
# The "Price" (D) and "Volume(1h)" (E) columns are refreshed with the latest
# quote snapshot. Two rows (43/44) also swap which coin ("dogwifhat" vs.
# "Cosmos") occupies that rank, so their Coin/Link/Price/Volume cells are
# rewritten together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "12.90", "0.999") would
# otherwise be silently re-typed as a numeric Value by Excel, which both
# changes the cell's stored type and can mangle the literal text (dropping a
# trailing zero, flipping "2.89" into a binary-float "2.8900000000000001",
# etc.). Force those through as text - restoring NumberFormat/Style right
# after keeps every cell's formatting identical to before the write.
function Set-TextCell {
    param($Sheet, $Addr, $Text)

    $cell = $Sheet.Range($Addr)
    $looksNumeric = $Text -match '^[+-]?\d+(\.\d+)*$'

    if ($looksNumeric) {
        $savedFormat = $cell.NumberFormat
        $cell.NumberFormat = '@'
        $cell.Value = $Text
        $cell.NumberFormat = $savedFormat
        $cell.Style = 'Normal'
    } else {
        $cell.Value = $Text
    }
}

$updates = @(
    @{ Cell = 'D2'; Value = '66.907.34' },
    @{ Cell = 'E2'; Value = '  +2.87%  ' },
    @{ Cell = 'D3'; Value = '3.102.16' },
    @{ Cell = 'E3'; Value = '  +5.48%  ' },
    @{ Cell = 'E4'; Value = '  -0.07%  ' },
    @{ Cell = 'D5'; Value = '581.23' },
    @{ Cell = 'E5'; Value = '  +2.51%  ' },
    @{ Cell = 'D6'; Value = '168.45' },
    @{ Cell = 'E6'; Value = '  +6.67%  ' },
    @{ Cell = 'D7'; Value = '0.999' },
    @{ Cell = 'E7'; Value = '  -0.11%  ' },
    @{ Cell = 'D8'; Value = '3.098.51' },
    @{ Cell = 'E8'; Value = '  +5.40%  ' },
    @{ Cell = 'E9'; Value = '  +1.42%  ' },
    @{ Cell = 'D10'; Value = '6.67' },
    @{ Cell = 'E10'; Value = '  -0.30%  ' },
    @{ Cell = 'E11'; Value = '  +3.09%  ' },
    @{ Cell = 'E12'; Value = '  +5.83%  ' },
    @{ Cell = 'E13'; Value = '  +2.52%  ' },
    @{ Cell = 'D14'; Value = '36.88' },
    @{ Cell = 'E14'; Value = '  +8.44%  ' },
    @{ Cell = 'E15'; Value = '  -0.71%  ' },
    @{ Cell = 'D16'; Value = '3.612.98' },
    @{ Cell = 'E16'; Value = '  +5.13%  ' },
    @{ Cell = 'D17'; Value = '66.946.33' },
    @{ Cell = 'E17'; Value = '  +2.60%  ' },
    @{ Cell = 'D18'; Value = '7.24' },
    @{ Cell = 'E18'; Value = '  +4.39%  ' },
    @{ Cell = 'D19'; Value = '3.100.67' },
    @{ Cell = 'E19'; Value = '  +5.19%  ' },
    @{ Cell = 'D20'; Value = '16.37' },
    @{ Cell = 'E20'; Value = '  +18.69%  ' },
    @{ Cell = 'D21'; Value = '470.26' },
    @{ Cell = 'E21'; Value = '  +5.55%  ' },
    @{ Cell = 'D22'; Value = '0.717' },
    @{ Cell = 'E22'; Value = '  +5.65%  ' },
    @{ Cell = 'E23'; Value = '  +4.90%  ' },
    @{ Cell = 'D24'; Value = '83.42' },
    @{ Cell = 'E24'; Value = '  +0.82%  ' },
    @{ Cell = 'D25'; Value = '2.37' },
    @{ Cell = 'E25'; Value = '  +9.09%  ' },
    @{ Cell = 'D26'; Value = '12.90' },
    @{ Cell = 'E26'; Value = '  +7.27%  ' },
    @{ Cell = 'D27'; Value = '10.18' },
    @{ Cell = 'E27'; Value = '  +2.16%  ' },
    @{ Cell = 'E28'; Value = '  +0.07%  ' },
    @{ Cell = 'D29'; Value = '8.08' },
    @{ Cell = 'E29'; Value = '  +2.04%  ' },
    @{ Cell = 'E30'; Value = '  +3.99%  ' },
    @{ Cell = 'E31'; Value = '  +4.63%  ' },
    @{ Cell = 'E32'; Value = '  +3.50%  ' },
    @{ Cell = 'D33'; Value = '28.32' },
    @{ Cell = 'E33'; Value = '  +4.16%  ' },
    @{ Cell = 'D34'; Value = '0.115' },
    @{ Cell = 'E34'; Value = '  +5.70%  ' },
    @{ Cell = 'E35'; Value = '  -0.08%  ' },
    @{ Cell = 'E36'; Value = '  +3.44%  ' },
    @{ Cell = 'D37'; Value = '5.92' },
    @{ Cell = 'E37'; Value = '  +4.46%  ' },
    @{ Cell = 'D38'; Value = '46.92' },
    @{ Cell = 'E38'; Value = '  +8.27%  ' },
    @{ Cell = 'D39'; Value = '2.09' },
    @{ Cell = 'E39'; Value = '  +6.03%  ' },
    @{ Cell = 'D40'; Value = '50.29' },
    @{ Cell = 'E41'; Value = '  +7.04%  ' },
    @{ Cell = 'D42'; Value = '0.124' },
    @{ Cell = 'E42'; Value = '  +4.29%  ' },
    @{ Cell = 'B43'; Value = 'Cosmos' },
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom' },
    @{ Cell = 'D43'; Value = '8.75' },
    @{ Cell = 'E43'; Value = '  +3.93%  ' },
    @{ Cell = 'B44'; Value = 'dogwifhat' },
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif' },
    @{ Cell = 'D44'; Value = '2.84' },
    @{ Cell = 'E44'; Value = '  +2.41%  ' },
    @{ Cell = 'D45'; Value = '392.03' },
    @{ Cell = 'E45'; Value = '  +2.68%  ' },
    @{ Cell = 'E46'; Value = '  +3.69%  ' },
    @{ Cell = 'D47'; Value = '2.762.40' },
    @{ Cell = 'E47'; Value = '  +1.10%  ' },
    @{ Cell = 'D48'; Value = '134.89' },
    @{ Cell = 'E48'; Value = '  +2.36%  ' },
    @{ Cell = 'E49'; Value = '  -0.01%  ' },
    @{ Cell = 'D50'; Value = '24.84' },
    @{ Cell = 'E50'; Value = '  +7.42%  ' },
    @{ Cell = 'E51'; Value = '  +4.70%  ' }
)

foreach ($u in $updates) {
    Set-TextCell $ws $u.Cell $u.Value
}
